$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.68"
$ws.Range("E2").Value = "'0.85%"
$ws.Range("D3").Value = "'26.88"
$ws.Range("E3").Value = "'-0.39%"
$ws.Range("D4").Value = "'4.644"
$ws.Range("E4").Value = "'-0.20%"
$ws.Range("D5").Value = "'0.05958"
$ws.Range("E5").Value = "'1.12%"
$ws.Range("D6").Value = "'6.639"
$ws.Range("E6").Value = "'0.07%"
$ws.Range("D7").Value = "'0.8561"
$ws.Range("E7").Value = "'-0.66%"
$ws.Range("D8").Value = "'0.9223"
$ws.Range("E8").Value = "'-1.16%"
$ws.Range("D9").Value = "'0.1384"
$ws.Range("E9").Value = "'-1.58%"
$ws.Range("D10").Value = "'0.04269"
$ws.Range("E10").Value = "'14.36%"
$ws.Range("D11").Value = "'0.07018"
$ws.Range("E11").Value = "'-1.01%"
$ws.Range("D12").Value = "'0.02977"
$ws.Range("E12").Value = "'-7.86%"
$ws.Range("D13").Value = "'0.09112"
$ws.Range("E13").Value = "'-1.14%"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'-1.12%"
$ws.Range("E15").Value = "'-0.01%"
$ws.Range("D16").Value = "'0.006086"
$ws.Range("E16").Value = "'0.14%"
$ws.Range("D17").Value = "'3.454"
$ws.Range("E17").Value = "'-1.77%"
$ws.Range("D18").Value = "'3.124"
$ws.Range("E18").Value = "'-2.09%"
$ws.Range("E19").Value = "'-2.19%"
$ws.Range("E20").Value = "'0.18%"
$ws.Range("D22").Value = "'4.013"
$ws.Range("E22").Value = "'4.37%"
$ws.Range("D23").Value = "'0.04216"
$ws.Range("E23").Value = "'-0.26%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-0.34%"
$ws.Range("D25").Value = "'0.004030"
$ws.Range("E25").Value = "'-5.80%"
$ws.Range("E27").Value = "'-11.63%"
$ws.Range("D40").Value = "'0.03822"
$ws.Range("E40").Value = "'-0.07%"
$ws.Range("D41").Value = "'0.1110"
$ws.Range("E41").Value = "'1.15%"
$ws.Range("D42").Value = "'0.003786"
$ws.Range("E42").Value = "'-39.02%"
$ws.Range("D43").Value = "'0.002427"
$ws.Range("E43").Value = "'10.36%"
$ws.Range("E44").Value = "'31.49%"
$ws.Range("D45").Value = "'0.00005154"
$ws.Range("E45").Value = "'-5.60%"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("E47").Value = "'-17.03%"
$ws.Range("D48").Value = "'0.2210"
$ws.Range("E48").Value = "'9,600.58%"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E50").Value = "'-0.07%"
